$d = $word.ActiveDocument

# Anchor: the (currently empty) paragraph immediately after the
# "3.1 User Interfaces" heading. This paragraph keeps its own pPr
# (paragraph-mark bold) but its empty run gets the first new sentence.
$find = $d.Content
$find.Find.Execute("3.1 User Interfaces", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara = $find.Paragraphs(1)
$p0 = $headingPara.Next()

# Insert five new paragraphs right after $p0 (before the existing
# "Color Palette..." paragraph). InsertParagraphBefore on a range
# collapsed at the start of a following paragraph creates a clean new
# paragraph (no stray empty <w:t/> runs) that inherits that paragraph's
# (non-bold) paragraph-mark formatting.
$insertionPoint = $p0.Next().Range.Duplicate
$insertionPoint.Collapse(1)

$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()
$insertionPoint.InsertParagraphBefore()

$p1 = $p0.Next()   # "At the bottom of the login page..."
$p2 = $p1.Next()   # empty
$p3 = $p2.Next()   # empty
$p4 = $p3.Next()   # "Upon successful login..."
$p5 = $p4.Next()   # empty

$p0.Range.InsertAfter("A first-time user of the app should see the Role Selection screen and depending on the selected role, the user will be redirected to the role-specific login page. Users can login using the email used to sign in or their google account (which will require NSU ID and a scanned picture). ")
$p1.Range.InsertAfter("At the bottom of the login page, there should be a Sign-Up button that redirects users to the Sign-Up webpage using the default browser.")
$p4.Range.InsertAfter("Upon successful login, the user should be at the main homepage where they have to fill the required fields to submit a complaint. The fields require information such as the complaint description, who the complaint is against, evidence submission, who should review the complaint.")

Write-Output "Text inserted."
